{"js": "// Replace each \"oldEquation\" cell text with its \"newEquation\" counterpart.\n// The document contains a 5-column table of three-digit \u00f7 one-digit division\n// problems; this diff swaps the operands/quotient/remainder text in 25 of\n// the populated cells while leaving everything else (fonts, sizes,\n// paragraph alignment, empty rows) untouched.\nconst replacements = [\n  [\"773\u00f77=110, 3\", \"186\u00f76=31, 0\"],\n  [\"127\u00f76=21, 1\", \"326\u00f79=36, 2\"],\n  [\"577\u00f79=64, 1\", \"960\u00f74=240, 0\"],\n  [\"366\u00f78=45, 6\", \"405\u00f74=101, 1\"],\n  [\"457\u00f73=152, 1\", \"627\u00f79=69, 6\"],\n  [\"554\u00f72=277, 0\", \"708\u00f76=118, 0\"],\n  [\"601\u00f76=100, 1\", \"474\u00f76=79, 0\"],\n  [\"492\u00f77=70, 2\", \"620\u00f78=77, 4\"],\n  [\"533\u00f75=106, 3\", \"489\u00f74=122, 1\"],\n  [\"688\u00f74=172, 0\", \"521\u00f79=57, 8\"],\n  [\"675\u00f75=135, 0\", \"304\u00f79=33, 7\"],\n  [\"779\u00f76=129, 5\", \"902\u00f73=300, 2\"],\n  [\"458\u00f72=229, 0\", \"251\u00f73=83, 2\"],\n  [\"336\u00f75=67, 1\", \"471\u00f74=117, 3\"],\n  [\"268\u00f76=44, 4\", \"845\u00f79=93, 8\"],\n  [\"858\u00f72=429, 0\", \"802\u00f78=100, 2\"],\n  [\"878\u00f77=125, 3\", \"515\u00f78=64, 3\"],\n  [\"701\u00f73=233, 2\", \"939\u00f72=469, 1\"],\n  [\"751\u00f72=375, 1\", \"201\u00f76=33, 3\"],\n  [\"849\u00f75=169, 4\", \"728\u00f77=104, 0\"],\n  [\"227\u00f78=28, 3\", \"298\u00f76=49, 4\"],\n  [\"162\u00f78=20, 2\", \"278\u00f78=34, 6\"],\n  [\"630\u00f72=315, 0\", \"576\u00f74=144, 0\"],\n  [\"804\u00f73=268, 0\", \"583\u00f79=64, 7\"],\n  [\"713\u00f79=79, 2\", \"404\u00f72=202, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"oldEquation\" cell text with its \"newEquation\" counterpart.\n# The document contains a 5-column table of three-digit division-by-one-digit\n# problems; this diff swaps the operands/quotient/remainder text in 25 of the\n# populated cells while leaving everything else (fonts, sizes, paragraph\n# alignment, empty rows) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"773\u00f77=110, 3\", \"186\u00f76=31, 0\"),\n    @(\"127\u00f76=21, 1\", \"326\u00f79=36, 2\"),\n    @(\"577\u00f79=64, 1\", \"960\u00f74=240, 0\"),\n    @(\"366\u00f78=45, 6\", \"405\u00f74=101, 1\"),\n    @(\"457\u00f73=152, 1\", \"627\u00f79=69, 6\"),\n    @(\"554\u00f72=277, 0\", \"708\u00f76=118, 0\"),\n    @(\"601\u00f76=100, 1\", \"474\u00f76=79, 0\"),\n    @(\"492\u00f77=70, 2\", \"620\u00f78=77, 4\"),\n    @(\"533\u00f75=106, 3\", \"489\u00f74=122, 1\"),\n    @(\"688\u00f74=172, 0\", \"521\u00f79=57, 8\"),\n    @(\"675\u00f75=135, 0\", \"304\u00f79=33, 7\"),\n    @(\"779\u00f76=129, 5\", \"902\u00f73=300, 2\"),\n    @(\"458\u00f72=229, 0\", \"251\u00f73=83, 2\"),\n    @(\"336\u00f75=67, 1\", \"471\u00f74=117, 3\"),\n    @(\"268\u00f76=44, 4\", \"845\u00f79=93, 8\"),\n    @(\"858\u00f72=429, 0\", \"802\u00f78=100, 2\"),\n    @(\"878\u00f77=125, 3\", \"515\u00f78=64, 3\"),\n    @(\"701\u00f73=233, 2\", \"939\u00f72=469, 1\"),\n    @(\"751\u00f72=375, 1\", \"201\u00f76=33, 3\"),\n    @(\"849\u00f75=169, 4\", \"728\u00f77=104, 0\"),\n    @(\"227\u00f78=28, 3\", \"298\u00f76=49, 4\"),\n    @(\"162\u00f78=20, 2\", \"278\u00f78=34, 6\"),\n    @(\"630\u00f72=315, 0\", \"576\u00f74=144, 0\"),\n    @(\"804\u00f73=268, 0\", \"583\u00f79=64, 7\"),\n    @(\"713\u00f79=79, 2\", \"404\u00f72=202, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
